# Progress.xlsx update — fill in missing/placeholder lab scores.
$wb = $excel.ActiveWorkbook

# ---- Sheet "БИВТ-22-17" (1st sheet) ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G2").Value = 5
$ws1.Range("F5").Value = 5
$ws1.Range("F9").Value = 5
$ws1.Range("F11").Value = 5
$ws1.Range("C15").Value = 5
$ws1.Range("F17").Value = 5
$ws1.Range("F19").Value = 5
$ws1.Range("J23").Value = 5
$ws1.Range("G25").Value = 5
$ws1.Range("G28").Value = 5

# ---- Sheet "БИВТ-22-18" (2nd sheet) ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("H3").Value = 5
$ws2.Range("K7").Value = 5
$ws2.Range("E8").Value = "pass"
$ws2.Range("F11").Value = "failed 3.5"
$ws2.Range("E12").Value = "wait al12"
$ws2.Range("E13").Value = 5
$ws2.Range("E14").Value = 5
$ws2.Range("D15").Value = 5
$ws2.Range("D16").Value = 5
$ws2.Range("G23").Value = 5

# ---- Selections / active sheet to mirror the author's final view ----
$ws1.Range("F5").Select()

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B31").Select()

$ws2.Activate()
$excel.ActiveWindow.Zoom = 145
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("F12").Select()
